$d = $word.ActiveDocument

# Replacement 1
$old1 = "G. H. Meeten and A. N. North, " + [string]([char]0x201C) + "Refractive index measurement of absorbing and turbid fluids by reflection near the critical angle," + [string]([char]0x201D) + " Meas. Sci. Technol. 6(2), 214" + [string]([char]0x2013) + "221 (1995). "
$new1 = "G. H. Meeten and A. N. North, " + [string]([char]0x201C) + "Refractive index measurement of absorbing and turbid fluids by reflection near the critical angle" + [string]([char]0x201D) + " Meas. Sci. Technol. 6(2), 214" + [string]([char]0x2013) + "221 (1995), DOI: 10.1088/0957-0233/6/2/014."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# Replacement 2
$old2 = "A. J. J" + [string]([char]0x00E4) + [string]([char]0x00E4) + "skel" + [string]([char]0x00E4) + "inen, K. E. Peiponen, and J. A. R" + [string]([char]0x00E4) + "ty, " + [string]([char]0x201C) + "On reflectometric measurement of a refractive index of milk," + [string]([char]0x201D) + " J. Dairy Sci. 84(1), 38" + [string]([char]0x2013) + "43 (2001). "
$new2 = "A. J. J" + [string]([char]0x00E4) + [string]([char]0x00E4) + "skel" + [string]([char]0x00E4) + "inen, K. E. Peiponen, and J. A. R" + [string]([char]0x00E4) + "ty, " + [string]([char]0x201C) + "On reflectometric measurement of a refractive index of milk" + [string]([char]0x201D) + " J. Dairy Sci. 84(1), 38" + [string]([char]0x2013) + "43 (2001), DOI: doi.org/10.3168/jds.S0022-0302(01)74449-9."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# Replacement 3
$old3 = "W. R. Calhoun, H. Maeta, A. Combs, L. M. Bali, and S. Bali, " + [string]([char]0x201C) + "Measurement of the refractive index of highly turbid media," + [string]([char]0x201D) + " Opt. Lett. 35(8), 1224" + [string]([char]0x2013) + "1226 (2010)."
$new3 = "W. R. Calhoun, H. Maeta, A. Combs, L. M. Bali, and S. Bali, " + [string]([char]0x201C) + "Measurement of the refractive index of highly turbid media" + [string]([char]0x201D) + " Opt. Lett. 35(8), 1224" + [string]([char]0x2013) + "1226 (2010)."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

